$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Contest 38 (row 47, "KXI vs DC") results entered.
# ------------------------------------------------------------------
$ws.Range("E47").Value = 20
$ws.Range("H47").Value = 80
$ws.Range("K47").Value = 60
$ws.Range("N47").Value = 40
$ws.Range("Q47").Value = 100
$ws.Range("T47").Value = 0

# ------------------------------------------------------------------
# Row 56 was the blank "next contest" template row. Insert a new blank
# row right there first: the old row 56 (with its template formulas
# still referencing the literal rank "2") slides down to become row 57,
# which is exactly the new blank template row for the contest after
# this one. The freshly inserted row 56 is what we fill in as contest
# 46 ("SRH vs DC").
# ------------------------------------------------------------------
$ws.Rows("56:56").Insert()

# Copy the (now row 57) template's formatting into the new row 56 so it
# keeps matching the rest of the results table.
$ws.Range("A57:C57").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("D57:E57").Copy()
$ws.Range("D56").PasteSpecial(-4122)
$ws.Range("G57:H57").Copy()
$ws.Range("G56").PasteSpecial(-4122)
$ws.Range("J57:K57").Copy()
$ws.Range("J56").PasteSpecial(-4122)
$ws.Range("M57:N57").Copy()
$ws.Range("M56").PasteSpecial(-4122)
$ws.Range("P57:Q57").Copy()
$ws.Range("P56").PasteSpecial(-4122)
$ws.Range("S57:T57").Copy()
$ws.Range("S56").PasteSpecial(-4122)

$ws.Range("A56").Value = 46
$ws.Range("B56").Value = 2
$ws.Range("C56").Value = "SRH vs DC"

$ws.Range("D56").Formula = '=IF(ISERROR(VLOOKUP(RANK(E56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(E56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("G56").Formula = '=IF(ISERROR(VLOOKUP(RANK(H56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(H56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("J56").Formula = '=IF(ISERROR(VLOOKUP(RANK(K56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(K56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("M56").Formula = '=IF(ISERROR(VLOOKUP(RANK(N56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(N56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("P56").Formula = '=IF(ISERROR(VLOOKUP(RANK(Q56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(Q56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'
$ws.Range("S56").Formula = '=IF(ISERROR(VLOOKUP(RANK(T56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE)),"",VLOOKUP(RANK(T56, ($T56,$Q56,$N56,$K56,$H56,$E56), 0),  score, $B56+1, FALSE))'

# ------------------------------------------------------------------
# Conditional formatting ranges don't follow a plain row insert in this
# runtime, so move each rule's target down to the (shifted) totals row
# by hand.
# ------------------------------------------------------------------
$ws.Range("E60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E61"))
$ws.Range("H60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H61"))
$ws.Range("K60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("K61"))
$ws.Range("N60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("N61"))
$ws.Range("Q60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("Q61"))
$ws.Range("T60").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("T61"))

# Keep the remembered selection in sync with the shifted totals cell.
$ws.Range("U61").Select()
